# Generate Report for Handoff
# Mark the 156c5445-0461-4b37-9980-0985a72c26a4.md file as "Ready for handoff"
# across the Overview, zh-cn and de-de sheets, updating handoff timestamps.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row for 156c5445-...md (row 3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-03-22 10:18:54"

# --- zh-cn sheet: row for 156c5445-...md (row 3) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "2016-03-22 10:18:50"

# --- de-de sheet: row for 156c5445-...md (row 3) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "2016-03-22 10:18:54"
